# Update canonical URL base from http://example.org/ to http://example.org/ig/example/
# and refresh the Date metadata value, per commit "changed canonical to example.org/ig/example".

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "http://example.org/ig/example/ValueSet/gvhd-all-valueset"
$wsMetadata.Range("B8").Value = "2023-04-26T11:15:05-05:00"

$wsInclude1 = $wb.Worksheets.Item("Include ValueSets")
$wsInclude1.Range("A2").Value = "http://example.org/ig/example/ValueSet/gvhd-sct-codes"

$wsInclude2 = $wb.Worksheets.Item("Include ValueSets 2")
$wsInclude2.Range("A2").Value = "http://example.org/ig/example/ValueSet/gvhd-icd10-codes"
